$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of A18 (long description text) and A20 (file URLs text).
# These are the top-left cells of merged ranges A18:F18 and A20:F20.
$nbsp = [char]0x00A0
$descricaoText = "Na montagem das laterais foi observado um erro no componente: 240642 - GANCHO MENOR P/ F.RÁPIDO CÓD.FT 0040-C - Saldo: 300,00 - Consumo Médio" + $nbsp + "Mensal:" + $nbsp + "87,67. Onde foram inspecionando os 100,00 na produção mais os 300,00 no almoxarifado. Do total dos 400,00 foi encontrado 47,00 unidades com não conformidades que não podem ser aproveitadas. "
$arquivosText = "https://cemag.monday.com/protected_static/12861583/resources/947643825/ev%202.jpg, https://cemag.monday.com/protected_static/12861583/resources/947643833/ev%201.jpg"

$ws.Range("A18").Value = $arquivosText
$ws.Range("A20").Value = $descricaoText

# Clear the "last updated" columns (Pessoas / Status / Avaliação) for row 23,
# and clear the whole rows 24-26 in the Tratamento table.
# Use ClearContents (not Clear) so cell formatting/style stays intact.
$ws.Range("D23:F23").ClearContents()
$ws.Range("A24:F25").ClearContents()
$ws.Range("A26:F26").ClearContents()
